$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.66
$ws.Range("G2").Value = 1.67
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 4.4
$ws.Range("N2").Value = 4
$ws.Range("Q2").Value = 1.89
$ws.Range("T2").Value = 1.93
$ws.Range("W2").Value = 2.48
$ws.Range("AF2").Value = 9.6
$ws.Range("AI2").Value = 85
$ws.Range("AN2").Value = 9.6
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5.9
$ws.Range("H5").Value = 1.66
$ws.Range("I5").Value = 1.72
$ws.Range("K5").Value = 4.6
$ws.Range("L5").Value = 1.29
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 2.1
$ws.Range("Q5").Value = 1.74
$ws.Range("R5").Value = 1.43
$ws.Range("S5").Value = 2.84
$ws.Range("T5").Value = 1.76
$ws.Range("U5").Value = 2.06
$ws.Range("V5").Value = 2.36
$ws.Range("W5").Value = 1.2
$ws.Range("X5").Value = 18.5
$ws.Range("Y5").Value = 10.5
$ws.Range("Z5").Value = 11
$ws.Range("AA5").Value = 17.5
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AE5").Value = 17.5
$ws.Range("AG5").Value = 22
$ws.Range("AH5").Value = 1000
$ws.Range("AK5").Value = 130
$ws.Range("AL5").Value = 140
$ws.Range("AM5").Value = 180
$ws.Range("AN5").Value = 180
$ws.Range("AO5").Value = 1000
$ws.Range("F6").Value = 1.42
$ws.Range("H6").Value = 8
$ws.Range("I6").Value = 11.5
$ws.Range("J6").Value = 4.6
$ws.Range("K6").Value = 5.5
$ws.Range("L6").Value = 1.22
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 4.3
$ws.Range("O6").Value = 1.21
$ws.Range("P6").Value = 2.14
$ws.Range("Q6").Value = 1.72
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 2.78
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.86
$ws.Range("V6").Value = 1.09
$ws.Range("W6").Value = 2.92
$ws.Range("X6").Value = 24
$ws.Range("Y6").Value = 36
$ws.Range("AB6").Value = 10.5
$ws.Range("AC6").Value = 14
$ws.Range("AD6").Value = 40
$ws.Range("AF6").Value = 11
$ws.Range("AG6").Value = 12.5
$ws.Range("AH6").Value = 32
$ws.Range("AJ6").Value = 15
$ws.Range("AK6").Value = 19
$ws.Range("AL6").Value = 44
$ws.Range("AN6").Value = 8
$ws.Range("F7").Value = 1.82
$ws.Range("G7").Value = 1.89
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 5.1
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 1.4
$ws.Range("N7").Value = 3.65
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 1.89
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.33
$ws.Range("S7").Value = 3.4
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 1.25
$ws.Range("W7").Value = 2.12
$ws.Range("X7").Value = 18.5
$ws.Range("AG7").Value = 10.5
$ws.Range("G9").Value = 2.52
$ws.Range("H9").Value = 3.25
$ws.Range("K9").Value = 3.65
$ws.Range("W9").Value = 1.66
$ws.Range("AE9").Value = 48
$ws.Range("F10").Value = 3.9
$ws.Range("G10").Value = 4.4
$ws.Range("H10").Value = 2.1
$ws.Range("I10").Value = 2.22
$ws.Range("J10").Value = 3.35
$ws.Range("K10").Value = 3.65
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 3.3
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 1.77
$ws.Range("Q10").Value = 1.95
$ws.Range("R10").Value = 1.29
$ws.Range("S10").Value = 3.85
$ws.Range("T10").Value = 1.84
$ws.Range("U10").Value = 1.98
$ws.Range("V10").Value = 1.83
$ws.Range("X10").Value = 15
$ws.Range("Y10").Value = 10.5
$ws.Range("Z10").Value = 15.5
$ws.Range("AB10").Value = 16
$ws.Range("AC10").Value = 9.199999999999999
$ws.Range("AD10").Value = 13
$ws.Range("AG10").Value = 20
$ws.Range("AH10").Value = 24
$ws.Range("AO10").Value = 24
$ws.Range("N11").Value = 3.35
$ws.Range("O11").Value = 1.39
$ws.Range("Q11").Value = 2.16
$ws.Range("U11").Value = 1.76
$ws.Range("H12").Value = 2.88
$ws.Range("I12").Value = 2.92
$ws.Range("J12").Value = 3.3
$ws.Range("K12").Value = 3.35
$ws.Range("V12").Value = 1.52
$ws.Range("Z12").Value = 17
$ws.Range("AH12").Value = 19.5
$ws.Range("AM12").Value = 130
$ws.Range("N13").Value = 3.1
$ws.Range("V13").Value = 1.93
$ws.Range("G14").Value = 3.65
$ws.Range("K14").Value = 3.25
$ws.Range("N14").Value = 2.68
$ws.Range("P14").Value = 1.59
$ws.Range("W14").Value = 1.39
$ws.Range("G15").Value = 2.18
$ws.Range("I15").Value = 4.5
$ws.Range("V15").Value = 1.3
$ws.Range("W15").Value = 1.84
$ws.Range("F16").Value = 2.3
$ws.Range("U16").Value = 2.14
$ws.Range("H17").Value = 5.6
$ws.Range("I17").Value = 5.7
$ws.Range("W17").Value = 2.26
$ws.Range("AH17").Value = 21
$ws.Range("J18").Value = 3.4
$ws.Range("N18").Value = 3.3
$ws.Range("O18").Value = 1.39
$ws.Range("P18").Value = 1.79
$ws.Range("V18").Value = 1.81
$ws.Range("L19").Value = 1.35
$ws.Range("P19").Value = 2.14
$ws.Range("U20").Value = 1.67
$ws.Range("F21").Value = 1.48
$ws.Range("G21").Value = 1.63
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 9.199999999999999
$ws.Range("J21").Value = 3.85
$ws.Range("K21").Value = 5.1
$ws.Range("N21").Value = 2.82
$ws.Range("Q21").Value = 1.64
$ws.Range("R21").Value = 1.34
$ws.Range("S21").Value = 2.52
$ws.Range("V21").Value = 1.12
$ws.Range("W21").Value = 2.58
